$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'36.990.79"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = "'2.047.42"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'251.23"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').Value = "'0.668"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = "'58.58"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.28%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'61.01"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('D10').Value = "'0.385"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('D11').Value = "'0.0789"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.54%  '
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').Value = "'16.27"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.73%  '
$ws.Range('D14').Value = "'2.350.84"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').Value = "'0.804"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.09%  '
$ws.Range('D16').Value = "'5.58"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.37%  '
$ws.Range('D17').Value = "'2.051.08"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').Value = "'36.930.19"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').Value = "'16.76"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +15.23%  '
$ws.Range('D20').Value = "'75.04"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('D21').Value = "'0.0₃0908"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.24%  '
$ws.Range('D22').Value = "'5.43"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.40%  '
$ws.Range('D23').Value = "'237.39"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.62%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').Value = "'2.38"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.07%  '
$ws.Range('E26').Value = '  +11.53%  '
$ws.Range('D27').Value = "'169.11"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = "'9.28"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').Value = "'20.20"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.46%  '
$ws.Range('D30').Value = "'0.125"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('E31').Value = '  +5.53%  '
$ws.Range('D32').Value = "'4.73"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.84%  '
$ws.Range('D33').Value = "'0.0618"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('D34').Value = "'4.46"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.75%  '
$ws.Range('D35').Value = "'0.0889"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = "'2.26"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('E38').Value = '  -3.63%  '
$ws.Range('E39').Value = '  +18.52%  '
$ws.Range('D40').Value = "'1.35"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('D41').Value = "'17.75"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.47%  '
$ws.Range('D42').Value = "'0.0224"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('D43').Value = "'1.14"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.27%  '
$ws.Range('D44').Value = "'96.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.53%  '
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('D46').Value = "'4.63"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +14.10%  '
$ws.Range('D47').Value = "'2.46"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('D48').Value = "'1.284.01"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.93%  '
$ws.Range('E49').Value = '  -1.63%  '
$ws.Range('D50').Value = "'6.82"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.20%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = "'2.236.16"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.61%  '
